$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update test result values for row 6
$ws.Range("F6").Value = "Same as expected outcome."
$ws.Range("G6").Value = "Pass"

# Update the selected range/active cell
$ws.Range("G5:G6").Select()
